# Refresh the crypto price/volume table (D:Price, E:Volume 1h) with the
# latest scrape, per the GitHub Actions update job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2 = @("96.925.69", "  -0.38%  ")
  3 = @("3.681.32", "  +2.48%  ")
  4 = @($null, "  +0.08%  ")
  5 = @("239.75", "  -0.72%  ")
  6 = @($null, "  +11.32%  ")
  7 = @("657.62", "  +0.29%  ")
  8 = @($null, "  +0.41%  ")
  9 = @("1.09", "  +3.55%  ")
  10 = @($null, "  +0.09%  ")
  11 = @("3.679.33", "  +2.37%  ")
  12 = @("45.68", "  +2.57%  ")
  13 = @($null, "  +0.88%  ")
  14 = @($null, "  +5.79%  ")
  15 = @("4.367.17", "  +2.60%  ")
  16 = @($null, "  +4.30%  ")
  17 = @("96.701.02", "  -0.36%  ")
  18 = @("8.95", "  +9.99%  ")
  19 = @("3.679.99", "  +2.71%  ")
  20 = @("18.85", "  +4.46%  ")
  21 = @("12.76", "  +0.28%  ")
  22 = @($null, "  +1.55%  ")
  23 = @("533.08", "  +3.63%  ")
  24 = @("3.52", "  +0.58%  ")
  25 = @("7.17", "  +4.78%  ")
  26 = @($null, "  -0.54%  ")
  27 = @("102.51", "  +1.71%  ")
  28 = @("13.51", "  +3.84%  ")
  29 = @($null, "  +5.55%  ")
  30 = @("12.39", "  +4.98%  ")
  31 = @($null, "  +1.28%  ")
  32 = @($null, "  +0.01%  ")
  33 = @("1.92", "  +17.11%  ")
  34 = @($null, "  +1.82%  ")
  35 = @($null, "  +0.46%  ")
  36 = @("32.69", "  +3.16%  ")
  37 = @("659.66", "  +6.14%  ")
  38 = @("0.598", "  +5.77%  ")
  39 = @("8.84", "  +1.11%  ")
  40 = @("0.160", "  +3.82%  ")
  41 = @("2.01", "  +2.53%  ")
  42 = @("6.62", "  +10.46%  ")
  43 = @("0.964", "  +4.25%  ")
  44 = @("38.81", "  +17.11%  ")
  45 = @($null, "  +0.04%  ")
  46 = @("0.0461", "  +4.85%  ")
  47 = @("0.429", "  +7.78%  ")
  48 = @("2.33", "  +1.58%  ")
  49 = @("3.74", "  +5.55%  ")
  50 = @("23.67", "  +0.18%  ")
  51 = @("8.72", "  +2.94%  ")
}

foreach ($row in $updates.Keys) {
  $pair = $updates[$row]
  $d = $pair[0]
  $e = $pair[1]
  if ($null -ne $d) {
    # Numeric-looking price strings must stay text (matches the source
    # data, which stores Price as inline-string, not a number) - force
    # the cell to Text before assigning, then drop the temporary format
    # so the cell ends up unstyled, same as before the edit.
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $d
    $ws.Range("D$row").ClearFormats()
  }
  $ws.Range("E$row").Value = $e
}
